$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.160.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").Value = "'1.840.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "

# Row 4
$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'240.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "

# Row 6
$ws.Range("D6").Value = "'0.6855"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.70%  "

# Row 7
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
$ws.Range("D8").Value = "'0.2993"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.57%  "

# Row 9
$ws.Range("D9").Value = "'0.07404"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.43%  "

# Row 10
$ws.Range("E10").Value = "  -2.30%  "

# Row 11
$ws.Range("D11").Value = "'0.07648"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.16%  "

# Row 12
$ws.Range("D12").Value = "'1.851.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "

# Row 13
$ws.Range("D13").Value = "'5.047"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.17%  "

# Row 14
$ws.Range("D14").Value = "'0.6809"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "

# Row 15
$ws.Range("D15").Value = "'87.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.64%  "

# Row 16
$ws.Range("D16").Value = "'6.135"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.37%  "

# Row 17
$ws.Range("D17").Value = "'29.174.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.40%  "

# Row 18
$ws.Range("D18").Value = "'0.000008148"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.43%  "

# Row 19
$ws.Range("D19").Value = "'2.084.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.66%  "

# Row 20
$ws.Range("D20").Value = "'228.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.66%  "

# Row 21
$ws.Range("D21").Value = "'12.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "

# Row 22
$ws.Range("D22").Value = "'1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").Value = "'7.368"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.43%  "

# Row 24
$ws.Range("E24").Value = "  +0.22%  "

# Row 25
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "'0.1445"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.28%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'159.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("D27").Value = "'8.747"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.42%  "

# Row 28
$ws.Range("E28").Value = "  -1.57%  "

# Row 29
$ws.Range("E29").Value = "  -1.77%  "

# Row 30
$ws.Range("D30").Value = "'4.268"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.59%  "

# Row 31
$ws.Range("E31").Value = "  -1.38%  "

# Row 32
$ws.Range("D32").Value = "'1.196"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "

# Row 33
$ws.Range("D33").Value = "'0.05253"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.21%  "

# Row 34
$ws.Range("D34").Value = "'0.7567"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.54%  "

# Row 35
$ws.Range("E35").Value = "  -3.36%  "

# Row 36
$ws.Range("D36").Value = "'1.134"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.58%  "

# Row 37
$ws.Range("D37").Value = "'2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$ws.Range("D38").Value = "'1.297.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.74%  "

# Row 39
$ws.Range("D39").Value = "'0.01828"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.60%  "

# Row 40
$ws.Range("E40").Value = "  -0.06%  "

# Row 41
$ws.Range("D41").Value = "'0.9361"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.87%  "

# Row 42
$ws.Range("D42").Value = "'5.963"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.82%  "

# Row 43
$ws.Range("D43").Value = "'104.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "

# Row 44
$ws.Range("D44").Value = "'0.9998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "

# Row 45
$ws.Range("D45").Value = "'1.988.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "

# Row 46
$ws.Range("D46").Value = "'0.5194"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47
$ws.Range("D47").Value = "'64.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "

# Row 48
$ws.Range("D48").Value = "'9.527"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.46%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000122"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.51%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.767"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "

# Row 51
$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").Value = "'0.07470"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.62%  "
